# codeforIATI/codelists data refresh: the SectorGroup codelist now lists
# codeforiati:group-code ahead of codeforiati:category-name, so swap the
# contents of columns E (category-name) and F (group-code) -- header row
# included -- while leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$rangeE = $ws.Range("E$($firstRow):E$($lastRow)")
$rangeF = $ws.Range("F$($firstRow):F$($lastRow)")
$rangeTemp = $ws.Range("Z$($firstRow):Z$($lastRow)")

$xlPasteValues = -4163

# Stash E's current contents (values + types) in the scratch column, then
# shuffle F -> E and the stash -> F. Copy/PasteSpecial(values) round-trips
# through Excel's own clipboard so the text-vs-number type of every cell
# (and its style) is preserved exactly, unlike a plain .Value re-assignment
# which would coerce numeric-looking text like "110" into a real number.
$rangeE.Copy() | Out-Null
$rangeTemp.PasteSpecial($xlPasteValues) | Out-Null

$rangeF.Copy() | Out-Null
$rangeE.PasteSpecial($xlPasteValues) | Out-Null

$rangeTemp.Copy() | Out-Null
$rangeF.PasteSpecial($xlPasteValues) | Out-Null

$rangeTemp.ClearContents() | Out-Null
$excel.CutCopyMode = 0
